$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 1.22
$ws.Range("H4").Value = 7
$ws.Range("I4").Value = 12
$ws.Range("J4").Value = 1.57
$ws.Range("L4").Value = 9
$ws.Range("O4").Value = 1.13
$ws.Range("P4").Value = 6
$ws.Range("Q4").Value = 1.4
$ws.Range("R4").Value = 3
$ws.Range("S4").Value = 1.22
$ws.Range("T4").Value = 4
$ws.Range("U4").Value = 1.95
$ws.Range("V4").Value = 1.8
$ws.Range("W4").Value = 10
$ws.Range("X4").Value = 7.5
$ws.Range("Z4").Value = 8
$ws.Range("AB4").Value = 23
$ws.Range("AC4").Value = 21
$ws.Range("AE4").Value = 23
$ws.Range("AG4").Value = 251
$ws.Range("AJ4").Value = 29
$ws.Range("AK4").Value = 151
$ws.Range("AL4").Value = 67
$ws.Range("AN4").Value = 3.4
$ws.Range("AQ4").Value = 12
$ws.Range("AR4").Value = 34
$ws.Range("AS4").Value = 101
$ws.Range("AT4").Value = 4
$ws.Range("AU4").Value = 9.5
$ws.Range("AZ4").Value = 201
$ws.Range("BB4").Value = 251
$ws.Range("BC4").Value = 501

# Row 5
$ws.Range("G5").Value = 1.22
$ws.Range("H5").Value = 6.5
$ws.Range("N5").Value = 26
$ws.Range("U5").Value = 1.7
$ws.Range("V5").Value = 2.05

# Row 6
$ws.Range("G6").Value = 1.55
$ws.Range("I6").Value = 5.5
$ws.Range("Q6").Value = 1.53
$ws.Range("R6").Value = 2.5
$ws.Range("U6").Value = 1.62
$ws.Range("V6").Value = 2.2
$ws.Range("X6").Value = 9
$ws.Range("Z6").Value = 12
$ws.Range("AC6").Value = 17
$ws.Range("AE6").Value = 15
$ws.Range("AL6").Value = 41
$ws.Range("BC6").Value = 401

# Row 7
$ws.Range("G7").Value = 1.44
$ws.Range("H7").Value = 4.5
$ws.Range("I7").Value = 7.5
$ws.Range("J7").Value = 1.91
$ws.Range("U7").Value = 1.75
$ws.Range("V7").Value = 2
$ws.Range("Z7").Value = 10
$ws.Range("AG7").Value = 201
$ws.Range("AK7").Value = 81
$ws.Range("AX7").Value = 34
$ws.Range("BA7").Value = 126

# Row 9
$ws.Range("G9").Value = 3.6
$ws.Range("H9").Value = 3.5
$ws.Range("Q9").Value = 1.92
$ws.Range("R9").Value = 1.98
$ws.Range("S9").Value = 1.4
$ws.Range("T9").Value = 2.75
$ws.Range("U9").Value = 1.75
$ws.Range("V9").Value = 2
$ws.Range("AA9").Value = 29
$ws.Range("AC9").Value = 11
$ws.Range("AH9").Value = 8
$ws.Range("AL9").Value = 17
$ws.Range("AM9").Value = 26
$ws.Range("AQ9").Value = 67
$ws.Range("AT9").Value = 2.75
$ws.Range("BB9").Value = 151

# Row 10
$ws.Range("G10").Value = 6
$ws.Range("I10").Value = 1.55
$ws.Range("M10").Value = 1.04
$ws.Range("N10").Value = 13
$ws.Range("Q10").Value = 1.67
$ws.Range("R10").Value = 2.2
$ws.Range("Y10").Value = 19
$ws.Range("Z10").Value = 67
$ws.Range("AC10").Value = 13
$ws.Range("AD10").Value = 7.5
$ws.Range("AI10").Value = 8
$ws.Range("AO10").Value = 29
$ws.Range("AQ10").Value = 101
$ws.Range("AS10").Value = 201
$ws.Range("AW10").Value = 3.6

# Row 19
$ws.Range("G19").Value = 3.25
$ws.Range("I19").Value = 2.63
$ws.Range("L19").Value = 3.6
$ws.Range("M19").Value = 1.18
$ws.Range("N19").Value = 4.5
$ws.Range("X19").Value = 13
$ws.Range("AI19").Value = 11
$ws.Range("AK19").Value = 29
$ws.Range("AN19").Value = 4.75
$ws.Range("AX19").Value = 19

# Row 37
$ws.Range("G37").Value = 3.1
$ws.Range("I37").Value = 2.4
$ws.Range("N37").Value = 8.5

# Row 48
$ws.Range("N48").Value = 9

# Row 107
$ws.Range("H107").Value = 5.4
$ws.Range("I107").Value = 1.23
$ws.Range("J107").Value = 6.3
$ws.Range("K107").Value = 3
$ws.Range("L107").Value = 1.57
$ws.Range("S107").Value = 1.08
$ws.Range("T107").Value = 6.4
$ws.Range("U107").Value = 1.3
$ws.Range("V107").Value = 3.38
$ws.Range("W107").Value = 65
$ws.Range("X107").Value = 120
$ws.Range("Y107").Value = 29
$ws.Range("AB107").Value = 32
$ws.Range("AC107").Value = 45
$ws.Range("AD107").Value = 15
$ws.Range("AE107").Value = 13.5
$ws.Range("AG107").Value = 65
$ws.Range("AH107").Value = 18.5
$ws.Range("AI107").Value = 12
$ws.Range("AJ107").Value = 9.5
$ws.Range("AK107").Value = 11.25
$ws.Range("AL107").Value = 8.75
$ws.Range("AN107").Value = 12.5
$ws.Range("AO107").Value = 40
$ws.Range("AP107").Value = 23
$ws.Range("AR107").Value = 110
$ws.Range("AS107").Value = 110
$ws.Range("AT107").Value = 6.2
$ws.Range("AU107").Value = 6.7
$ws.Range("AV107").Value = 25
$ws.Range("AW107").Value = 4.4
$ws.Range("AX107").Value = 5.8
$ws.Range("AY107").Value = 9
$ws.Range("AZ107").Value = 11.25
$ws.Range("BA107").Value = 17.5
$ws.Range("BB107").Value = 50
$ws.Range("BC107").Value = 200
